$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.570.25'
$ws.Range("E2").Value = '  -1.76%  '

$ws.Range("D3").Value = '2.613.19'
$ws.Range("E3").Value = '  +0.01%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '531.60'
$ws.Range("E5").Value = '  -1.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.01'
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.564'
$ws.Range("E8").Value = '  -0.55%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.96'
$ws.Range("E9").Value = '  +7.72%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1000'
$ws.Range("E10").Value = '  -3.02%  '

$ws.Range("E11").Value = '  -0.80%  '

$ws.Range("E12").Value = '  +0.57%  '

$ws.Range("D13").Value = '3.071.52'
$ws.Range("E13").Value = '  +0.02%  '

$ws.Range("D14").Value = '58.456.80'
$ws.Range("E14").Value = '  -1.82%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.77'
$ws.Range("E15").Value = '  -0.26%  '

$ws.Range("D16").Value = '2.609.95'
$ws.Range("E16").Value = '  +0.38%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000132'
$ws.Range("E17").Value = '  -1.29%  '

$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '334.70'
$ws.Range("E18").Value = '  -2.42%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.37'
$ws.Range("E19").Value = '  -0.34%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.07'
$ws.Range("E20").Value = '  -1.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.18'
$ws.Range("E21").Value = '  -3.78%  '

$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.30'
$ws.Range("E23").Value = '  -2.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.411'
$ws.Range("E24").Value = '  +0.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  +0.07%  '

$ws.Range("E26").Value = '  -1.70%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.12'
$ws.Range("E27").Value = '  -2.64%  '

$ws.Range("E28").Value = '  -0.02%  '

$ws.Range("D29").Value = '0.0₃0726'
$ws.Range("E29").Value = '  -3.31%  '

$ws.Range("E30").Value = '  -3.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.78'
$ws.Range("E31").Value = '  -1.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '150.78'
$ws.Range("E32").Value = '  +0.27%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.64'
$ws.Range("E33").Value = '  -1.24%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.88'
$ws.Range("E34").Value = '  -3.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.10'
$ws.Range("E35").Value = '  -2.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.818'
$ws.Range("E36").Value = '  -2.45%  '

$ws.Range("B37").Value = 'Stacks'
$ws.Range("C37").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.42'
$ws.Range("E37").Value = '  -4.18%  '

$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.812'
$ws.Range("E38").Value = '  -3.52%  '

$ws.Range("E39").Value = '  -0.62%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '281.76'
$ws.Range("E40").Value = '  +1.60%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.68'
$ws.Range("E42").Value = '  -0.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.589'
$ws.Range("E43").Value = '  -2.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0529'
$ws.Range("E44").Value = '  +0.56%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0934'
$ws.Range("E45").Value = '  -2.47%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.78'
$ws.Range("E46").Value = '  +0.22%  '

$ws.Range("E47").Value = '  -0.52%  '

$ws.Range("D48").Value = '1.929.02'
$ws.Range("E48").Value = '  -1.79%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.46'
$ws.Range("E49").Value = '  -1.97%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.85'
$ws.Range("E50").Value = '  -4.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '111.26'
$ws.Range("E51").Value = '  -1.06%  '
